# Auto-generated Excel COM-interop script
# Applies the 2025-09-17 daily crime-count update across all affected sheets
# in the Chicago violent-crime workbook (Citywide Totals, By Neighborhood, and
# every individual neighborhood sheet whose 2025 (column L) / historical totals changed).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 12).Value = 4804
$ws.Cells.Item(3, 12).Value = 5179
$ws.Cells.Item(4, 2).Value = 1716
$ws.Cells.Item(4, 5).Value = 2056
$ws.Cells.Item(4, 12).Value = 1269
$ws.Cells.Item(6, 12).Value = 4380
$ws.Cells.Item(7, 2).Value = 23348
$ws.Cells.Item(7, 5).Value = 26061
$ws.Cells.Item(7, 12).Value = 15936

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(7, 12).Value = 522
$ws.Cells.Item(8, 12).Value = 1059
$ws.Cells.Item(11, 12).Value = 259
$ws.Cells.Item(15, 12).Value = 116
$ws.Cells.Item(19, 12).Value = 437
$ws.Cells.Item(20, 12).Value = 399
$ws.Cells.Item(25, 12).Value = 94
$ws.Cells.Item(27, 12).Value = 142
$ws.Cells.Item(29, 12).Value = 873
$ws.Cells.Item(31, 12).Value = 160
$ws.Cells.Item(32, 12).Value = 21
$ws.Cells.Item(33, 12).Value = 729
$ws.Cells.Item(36, 12).Value = 207
$ws.Cells.Item(37, 12).Value = 592
$ws.Cells.Item(42, 12).Value = 518
$ws.Cells.Item(46, 12).Value = 35
$ws.Cells.Item(47, 12).Value = 111
$ws.Cells.Item(48, 12).Value = 206
$ws.Cells.Item(51, 12).Value = 200
$ws.Cells.Item(52, 12).Value = 320
$ws.Cells.Item(54, 12).Value = 334
$ws.Cells.Item(60, 12).Value = 102
$ws.Cells.Item(63, 2).Value = 420
$ws.Cells.Item(63, 5).Value = 389
$ws.Cells.Item(63, 12).Value = 45
$ws.Cells.Item(65, 12).Value = 311
$ws.Cells.Item(67, 12).Value = 548
$ws.Cells.Item(72, 12).Value = 62
$ws.Cells.Item(73, 12).Value = 125
$ws.Cells.Item(78, 12).Value = 209
$ws.Cells.Item(79, 12).Value = 421
$ws.Cells.Item(80, 12).Value = 51
$ws.Cells.Item(82, 12).Value = 23
$ws.Cells.Item(83, 12).Value = 353
$ws.Cells.Item(84, 12).Value = 154
$ws.Cells.Item(85, 12).Value = 817
$ws.Cells.Item(87, 12).Value = 48
$ws.Cells.Item(88, 12).Value = 173
$ws.Cells.Item(89, 12).Value = 231
$ws.Cells.Item(90, 12).Value = 160
$ws.Cells.Item(91, 12).Value = 220
$ws.Cells.Item(92, 12).Value = 45
$ws.Cells.Item(94, 12).Value = 195
$ws.Cells.Item(95, 12).Value = 220
$ws.Cells.Item(96, 12).Value = 179
$ws.Cells.Item(97, 12).Value = 136
$ws.Cells.Item(98, 12).Value = 86
$ws.Cells.Item(99, 12).Value = 274
$ws.Cells.Item(101, 2).Value = 23348
$ws.Cells.Item(101, 5).Value = 26061
$ws.Cells.Item(101, 12).Value = 15936

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(3, 12).Value = 52
$ws.Cells.Item(7, 12).Value = 179

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 12).Value = 177
$ws.Cells.Item(3, 12).Value = 174
$ws.Cells.Item(7, 12).Value = 522

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(4, 12).Value = 22
$ws.Cells.Item(7, 12).Value = 259

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(3, 12).Value = 69
$ws.Cells.Item(6, 12).Value = 62
$ws.Cells.Item(7, 12).Value = 231

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(3, 12).Value = 333
$ws.Cells.Item(7, 12).Value = 817

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(4, 12).Value = 21
$ws.Cells.Item(7, 12).Value = 320

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(3, 12).Value = 356
$ws.Cells.Item(7, 12).Value = 1059

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(6, 12).Value = 82
$ws.Cells.Item(7, 12).Value = 353

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 12).Value = 199
$ws.Cells.Item(3, 12).Value = 250
$ws.Cells.Item(6, 12).Value = 221
$ws.Cells.Item(7, 12).Value = 729

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(2, 12).Value = 87
$ws.Cells.Item(7, 12).Value = 220

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 12).Value = 174
$ws.Cells.Item(3, 12).Value = 203
$ws.Cells.Item(4, 12).Value = 33
$ws.Cells.Item(6, 12).Value = 164
$ws.Cells.Item(7, 12).Value = 592

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(2, 12).Value = 111
$ws.Cells.Item(6, 12).Value = 83
$ws.Cells.Item(7, 12).Value = 311

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(3, 12).Value = 115
$ws.Cells.Item(7, 12).Value = 274

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(2, 12).Value = 63
$ws.Cells.Item(7, 12).Value = 160

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(2, 12).Value = 159
$ws.Cells.Item(7, 12).Value = 548

$ws = $wb.Worksheets.Item('South Deering')
$ws.Cells.Item(6, 12).Value = 44
$ws.Cells.Item(7, 12).Value = 154

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(2, 12).Value = 61
$ws.Cells.Item(3, 12).Value = 81
$ws.Cells.Item(7, 12).Value = 334

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 12).Value = 261
$ws.Cells.Item(3, 12).Value = 329
$ws.Cells.Item(7, 12).Value = 873

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(3, 12).Value = 52
$ws.Cells.Item(7, 12).Value = 206

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(2, 12).Value = 155
$ws.Cells.Item(7, 12).Value = 437

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 12).Value = 149
$ws.Cells.Item(7, 12).Value = 518

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(6, 12).Value = 63
$ws.Cells.Item(7, 12).Value = 209

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Cells.Item(3, 12).Value = 10
$ws.Cells.Item(7, 12).Value = 35

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(2, 12).Value = 77
$ws.Cells.Item(3, 12).Value = 98
$ws.Cells.Item(7, 12).Value = 220

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(2, 12).Value = 138
$ws.Cells.Item(3, 12).Value = 150
$ws.Cells.Item(6, 12).Value = 91
$ws.Cells.Item(7, 12).Value = 421

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(2, 12).Value = 123
$ws.Cells.Item(6, 12).Value = 107
$ws.Cells.Item(7, 12).Value = 399

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(3, 12).Value = 61
$ws.Cells.Item(4, 12).Value = 16
$ws.Cells.Item(6, 12).Value = 53
$ws.Cells.Item(7, 12).Value = 207

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(3, 12).Value = 45
$ws.Cells.Item(7, 12).Value = 195

$ws = $wb.Worksheets.Item('East Side')
$ws.Cells.Item(2, 12).Value = 34
$ws.Cells.Item(7, 12).Value = 94

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(2, 12).Value = 42
$ws.Cells.Item(7, 12).Value = 111

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(3, 12).Value = 37
$ws.Cells.Item(7, 12).Value = 116

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Cells.Item(2, 12).Value = 20
$ws.Cells.Item(7, 12).Value = 86

$ws = $wb.Worksheets.Item('North Center')
$ws.Cells.Item(3, 12).Value = 9
$ws.Cells.Item(6, 12).Value = 13

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(2, 12).Value = 46
$ws.Cells.Item(7, 12).Value = 125

$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(2, 12).Value = 31
$ws.Cells.Item(7, 12).Value = 136

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Cells.Item(6, 12).Value = 17
$ws.Cells.Item(7, 12).Value = 45

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(2, 12).Value = 49
$ws.Cells.Item(7, 12).Value = 173

$ws = $wb.Worksheets.Item('Galewood')
$ws.Cells.Item(4, 12).Value = 3
$ws.Cells.Item(7, 12).Value = 21

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(4, 12).Value = 19
$ws.Cells.Item(7, 12).Value = 142

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(2, 12).Value = 53
$ws.Cells.Item(7, 12).Value = 160

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(2, 12).Value = 57
$ws.Cells.Item(3, 12).Value = 65
$ws.Cells.Item(7, 12).Value = 200

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Cells.Item(2, 12).Value = 33
$ws.Cells.Item(3, 12).Value = 36
$ws.Cells.Item(7, 12).Value = 102

$ws = $wb.Worksheets.Item('Old Town')
$ws.Cells.Item(3, 12).Value = 16
$ws.Cells.Item(7, 12).Value = 62

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Cells.Item(4, 12).Value = 5
$ws.Cells.Item(7, 12).Value = 23

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Cells.Item(2, 12).Value = 10
$ws.Cells.Item(7, 12).Value = 51

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Cells.Item(6, 12).Value = 18
$ws.Cells.Item(7, 12).Value = 48
